# Weekly update: a new daily price record was inserted for
# "Femacal de La Calera - Ciboulette" at what is now row 26,
# pushing the existing rows 26-182 down by one (to 27-183).
# The new row reuses the same Volumen/Precio/Unidad/Origen/etc.
# values that used to sit in row 26, only the Fecha (column D)
# is new (Excel serial date 44473).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$insertRow = 26

# Insert a blank row above the current row 26; this shifts the
# old row 26 (and everything below it) down to row 27, etc.
$ws.Rows.Item($insertRow).Insert()

# The row that used to be 26 is now at row 27 - copy its values
# into the newly-inserted row 26, then overwrite the date.
$sourceRow = $insertRow + 1
$columns = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R")

foreach ($col in $columns) {
    $srcCell = $ws.Range($col + $sourceRow)
    $dstCell = $ws.Range($col + $insertRow)
    $dstCell.Value = $srcCell.Value()
}

# New record's date (Fecha) - Excel serial date 44473 (2021-10-04)
$ws.Range("D" + $insertRow).Value = 44473
